# Apply the symbol-list refresh for the crypto tracker sheet (GitHub Actions scrape).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value.
$cellValues = [ordered]@{
    "D2" = "243.87"
    "D3" = "25.19"
    "B4" = "HuobiToken"
    "C4" = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
    "D4" = "5.169"
    "E4" = "3HuobiTokenHT"
    "B5" = "Cronos"
    "C5" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
    "D5" = "0.05735"
    "E5" = "4CronosCRO"
    "B6" = "KuCoinToken"
    "C6" = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
    "D6" = "6.499"
    "E6" = "5KuCoinTokenKCS"
    "B7" = "GateToken"
    "C7" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "D7" = "3.109"
    "E7" = "6GateTokenGT"
    "B8" = "MXToken"
    "C8" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "D8" = "0.8095"
    "E8" = "7MXTokenMX"
    "B9" = "FTXToken"
    "C9" = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
    "D9" = "0.8446"
    "E9" = "8FTXTokenFTT"
    "B10" = "WazirX"
    "C10" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D10" = "0.1338"
    "E10" = "9WazirXWRX"
    "B11" = "MandalaExchangeToken"
    "C11" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "D11" = "0.06953"
    "E11" = "10MandalaExchangeTokenMDX"
    "B12" = "BitrueCoin"
    "C12" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "D12" = "0.02829"
    "E12" = "11BitrueCoinBTR"
    "B13" = "BitMartToken"
    "C13" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D13" = "0.09369"
    "E13" = "12BitMartTokenBMX"
    "B14" = "BitForexToken"
    "C14" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D14" = "0.001511"
    "E14" = "13BitForexTokenBF"
    "B15" = "One"
    "C15" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "D15" = "0.0005997"
    "E15" = "14OneONE"
    "B16" = "TigerCash"
    "C16" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "D16" = "0.006085"
    "E16" = "15TigerCashTCH"
    "B17" = "LEO"
    "C17" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D17" = "3.500"
    "E17" = "16LEOLEO"
    "D18" = "2.053"
    "D19" = "0.3199"
    "D20" = "0.03126"
    "D21" = "0.1300"
    "D22" = "3.741"
    "D23" = "0.04665"
    "D24" = "0.1328"
    "D25" = "0.001236"
    "D26" = "0.004266"
    "D27" = "0.00009696"
    "E27" = "26NitroExNTXBestin24h"
    "D28" = "0.0001500"
    "D40" = "0.03609"
    "B41" = "BKEXToken"
    "C41" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "D41" = "0.1051"
    "E41" = "40BKEXTokenBKK"
    "B42" = "CEJI"
    "C42" = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
    "D42" = "0.002901"
    "E42" = "41CEJICEJI"
    "B43" = "KickToken"
    "C43" = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
    "D43" = "0.003164"
    "E43" = "42KickTokenKICK"
    "D44" = "0.007347"
    "D45" = "0.00005295"
    "D47" = "0.1499"
    "D48" = "0.002310"
    "D49" = "0.00002099"
    "D50" = "0.0001999"
}

# All of these values are numeric-looking strings (e.g. "243.87", "0.1300", "3.500") that
# must stay as text (matching the original inline-string cells) rather than being parsed
# into numbers and losing their formatting/trailing zeros. Force Text format first, assign
# the values, then restore the default "Normal" style so no stray number-format override
# remains on the cell.
foreach ($ref in $cellValues.Keys) {
    $ws.Range($ref).NumberFormat = "@"
}
foreach ($ref in $cellValues.Keys) {
    $ws.Range($ref).Value = $cellValues[$ref]
}
foreach ($ref in $cellValues.Keys) {
    $ws.Range($ref).Style = "Normal"
}
